$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3464964993005633, 1.65323645889881, 16.98373111632243, 6.48142807727062, 25.46489215179242)
    3  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    4  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    5  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    6  = @(0.3464964993005633, 0.3375848360084654, 0.1529057820181812, 6.48142807727062, 7.31841519459783)
    7  = @(0.1554434735375247, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.420375733316612)
    8  = @(3.182878228561681, 1.65323645889881, 157.8057217802531, 6.48142807727062, 169.1232645449842)
    9  = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
    10 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.964442013611383)
    11 = @(3.182878228561681, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 13.08449056854435)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
